# Fix the typo "biallic" -> "biallelic" in the "include all polymorphic
# sites (biallic AND multiallelic)" textbox, collapsing the three runs
# that made up the sentence back into a single run (matching the
# formatting already used by the run that contains the opening text).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("TextBox 260")
$tr = $shape.TextFrame.TextRange

# First, retype the misspelled word plus everything after it …
$tail = $tr.Characters(32, 100)
$tail.Text = "biallelic AND multiallelic)   "

# … then retype the whole line so it collapses into a single run that
# takes on the formatting of the first (correctly-tagged) run.
$whole = $tr.Characters(1, 200)
$whole.Text = "include all polymorphic sites (biallelic AND multiallelic)   "
